$d = $word.ActiveDocument

# 1. Update version string
$d.Content.Find.Execute("Verze: 1.00", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Verze: 1.01", 2)
